# Refresh the crypto symbol/price table (GitHub Actions daily scrape update).
# Only the Price (D) / Volume(1h) (E) columns move for most rows; rows 16-24
# additionally shift up because a new "HotbitToken" listing was inserted at
# row 16 (pushing the former LEO..BitKan block down one row and dropping the
# previous HotbitToken row that used to sit at row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: row number -> Coin, Link, Price, Volume(1h).
# $null means "leave this column alone".
$updates = @(
    @{ Row = 2;  D = "286.50";       E = "-9.95%" }
    @{ Row = 3;  D = "40.21";        E = "-2.51%" }
    @{ Row = 4;  D = "5.027";        E = "-4.06%" }
    @{ Row = 5;  D = "0.07286";      E = "-5.82%" }
    @{ Row = 6;  D = "4.288";        E = "-0.56%" }
    @{ Row = 7;  D = "1.527";        E = "-10.60%" }
    @{ Row = 8;  D = "0.9180";       E = "-3.55%" }
    @{ Row = 9;  D = "0.1198";       E = "-5.47%" }
    @{ Row = 10; D = "0.1709";       E = "-6.44%" }
    @{ Row = 11; D = "0.08640";      E = "-6.12%" }
    @{ Row = 12; D = "0.04170";      E = "-3.47%" }
    @{ Row = 13; D = "0.1050";       E = "-0.62%" }
    @{ Row = 14; D = "0.001274";     E = "0.31%" }
    @{ Row = 15; D = "0.005869";     E = "-0.21%" }

    @{ Row = 16; B = "HotbitToken";             C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";             D = "0.003780";    E = "-8.23%" }
    @{ Row = 17; B = "LEO";                     C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                     D = "3.399";       E = "1.30%" }
    @{ Row = 18; B = "BTSEToken";                C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";             D = "2.397";       E = "-1.16%" }
    @{ Row = 19; B = "BitpandaEcosystemToken";  C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";     D = "0.3282";      E = "-2.11%" }
    @{ Row = 20; B = "MCDex";                   C = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb";                     D = "7.847";       E = "4.30%" }
    @{ Row = 21; B = "ProBitToken";             C = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob";               D = "0.1344";      E = "-0.57%" }
    @{ Row = 22; B = "ZBToken";                 C = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb";                D = "0.2887";      E = "2.59%" }
    @{ Row = 23; B = "CoinExToken";             C = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet";           D = "0.03846";     E = "-4.60%" }
    @{ Row = 24; B = "BitKan";                  C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";                D = "0.001271";    E = "0.58%" }

    @{ Row = 25; D = "0.0001282";    E = "1.00%" }
    @{ Row = 26; D = "0.0003731";    E = "-95.03%" }

    @{ Row = 38; D = "0.02309";       E = "-9.08%" }
    @{ Row = 39; D = "0.04989";       E = "-6.83%" }
    @{ Row = 40; D = "0.006769";      E = "241.95%" }
    @{ Row = 41; D = "0.007688";      E = "-1.49%" }
    @{ Row = 42; D = "0.1266";        E = "-4.07%" }
    @{ Row = 43; D = "0.007376";      E = "0.42%" }
    @{ Row = 44; D = "0.007422";      E = "-2.01%" }
    @{ Row = 45; D = "0.3091";        E = "-10.00%" }
    @{ Row = 46; D = "0.00006451";    E = "-3.98%" }
    @{ Row = 47; D = "0.00000000752"; E = "0.32%" }
    @{ Row = 48; E = "14.13%" }
    @{ Row = 49; E = "0.05%" }
    @{ Row = 50; D = "0.00002105";    E = "0.32%" }
    @{ Row = 51; D = "0.0002005";     E = "0.32%" }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            # Briefly force Text format so values like "286.50" / "-9.95%" are
            # stored verbatim instead of being parsed into numbers/percents,
            # then restore the plain "Normal" style so we don't leave a stray
            # format/style change behind (the source diff only touches text).
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
            $cell.NumberFormat = "General"
            $cell.Style = "Normal"
        }
    }
}
